# 7.62x39 damage increased from 40 to 45
# "gun power" column (H) for ammo_7.62x39_fmj (row 19) and ammo_7.62x39_ap (row 20)
# goes from 1.05 to 1.17. Downstream columns J, K, E are formulas and recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("H19").Value = 1.17
$ws.Range("H20").Value = 1.17

$ws.Range("K13").Select()
